$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 14 (pushes existing rows 14..123 down to 15..124)
$ws.Rows.Item(14).Insert()

# Populate the newly inserted row 14 with the new weekly price report entry
$ws.Range("A14").Value = 10
$ws.Range("B14").Value = "Vega Modelo de Temuco"
$ws.Range("C14").Value = "La Araucanía"
$ws.Range("D14").Value = 44532
$ws.Range("E14").Value = 9
$ws.Range("F14").Value = 100114007
$ws.Range("G14").Value = "Jengibre"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 30
$ws.Range("K14").Value = 20000
$ws.Range("L14").Value = 20000
$ws.Range("M14").Value = 20000
$ws.Range("N14").Value = "$/caja 13 kilos"
$ws.Range("O14").Value = "Perú"
$ws.Range("P14").Value = 1538
$ws.Range("Q14").Value = 13
$ws.Range("R14").Value = "Hortaliza"
